$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format while we write values, so that numeric-looking
# strings (e.g. "236.66") are not auto-converted into real numbers by Excel.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '29.186.04'
$ws.Range("E2").Value = '  -0.55%  '
$ws.Range("D3").Value = '1.827.50'
$ws.Range("E3").Value = '  -0.74%  '
$ws.Range("D4").Value = '0.9988'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '236.66'
$ws.Range("E5").Value = '  -1.33%  '
$ws.Range("D6").Value = '0.6072'
$ws.Range("E6").Value = '  -3.48%  '
$ws.Range("D7").Value = '0.9996'
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '0.07089'
$ws.Range("E8").Value = '  -4.87%  '
$ws.Range("D9").Value = '0.2812'
$ws.Range("E9").Value = '  -2.87%  '
$ws.Range("D10").Value = '23.79'
$ws.Range("E10").Value = '  -4.96%  '
$ws.Range("D11").Value = '0.07660'
$ws.Range("E11").Value = '  -0.84%  '
$ws.Range("D12").Value = '1.809.90'
$ws.Range("E12").Value = '  -4.56%  '
$ws.Range("D13").Value = '4.820'
$ws.Range("E13").Value = '  -3.08%  '
$ws.Range("D14").Value = '0.00001005'
$ws.Range("E14").Value = '  -2.75%  '
$ws.Range("D15").Value = '0.6348'
$ws.Range("E15").Value = '  -6.16%  '
$ws.Range("D16").Value = '2.068.57'
$ws.Range("E16").Value = '  -0.30%  '
$ws.Range("D17").Value = '79.16'
$ws.Range("E17").Value = '  -3.26%  '
$ws.Range("D18").Value = '5.900'
$ws.Range("E18").Value = '  -5.49%  '
$ws.Range("D19").Value = '29.203.22'
$ws.Range("E19").Value = '  -0.47%  '
$ws.Range("D20").Value = '227.96'
$ws.Range("E20").Value = '  -0.44%  '
$ws.Range("D21").Value = '11.79'
$ws.Range("E21").Value = '  -4.33%  '
$ws.Range("D22").Value = '0.9996'
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("D23").Value = '7.022'
$ws.Range("E23").Value = '  -4.75%  '
$ws.Range("D24").Value = '0.9993'
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("D25").Value = '154.46'
$ws.Range("E25").Value = '  -2.24%  '
$ws.Range("D26").Value = '8.069'
$ws.Range("E26").Value = '  -5.42%  '
$ws.Range("D27").Value = '0.1298'
$ws.Range("E27").Value = '  -3.74%  '
$ws.Range("D28").Value = '16.56'
$ws.Range("E28").Value = '  -4.97%  '
$ws.Range("D29").Value = '1.478'
$ws.Range("E29").Value = '  +1.42%  '
$ws.Range("D30").Value = '0.06462'
$ws.Range("E30").Value = '  -5.76%  '
$ws.Range("E31").Value = '  -1.95%  '
$ws.Range("D32").Value = '3.827'
$ws.Range("E32").Value = '  -5.54%  '
$ws.Range("D33").Value = '3.806'
$ws.Range("E33").Value = '  -6.38%  '
$ws.Range("D34").Value = '1.126'
$ws.Range("E34").Value = '  -1.17%  '
$ws.Range("D35").Value = '1.747'
$ws.Range("E35").Value = '  -4.48%  '
$ws.Range("D36").Value = '0.6495'
$ws.Range("E36").Value = '  -7.15%  '
$ws.Range("D37").Value = '2.547'
$ws.Range("E37").Value = '  -1.27%  '
$ws.Range("D38").Value = '2.750'
$ws.Range("E38").Value = '  -2.53%  '
$ws.Range("D39").Value = '1.213.93'
$ws.Range("E39").Value = '  -1.95%  '
$ws.Range("D40").Value = '0.01749'
$ws.Range("E40").Value = '  -5.26%  '
$ws.Range("D41").Value = '6.494'
$ws.Range("E41").Value = '  -4.57%  '
$ws.Range("D42").Value = '0.9338'
$ws.Range("E42").Value = '  -0.74%  '
$ws.Range("D43").Value = '0.9989'
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("D44").Value = '101.01'
$ws.Range("E44").Value = '  -0.08%  '
$ws.Range("D45").Value = '1.983.47'
$ws.Range("E45").Value = '  -0.29%  '
$ws.Range("D46").Value = '63.05'
$ws.Range("E46").Value = '  -3.55%  '
$ws.Range("D47").Value = '0.00000000119'
$ws.Range("E47").Value = '  -0.20%  '
$ws.Range("D48").Value = '1.611'
$ws.Range("D49").Value = '8.588'
$ws.Range("E49").Value = '  -4.08%  '
$ws.Range("D50").Value = '0.1076'
$ws.Range("E50").Value = '  -5.85%  '
$ws.Range("D51").Value = '0.05526'
$ws.Range("E51").Value = '  -2.69%  '

# Restore the default (Normal) style on column D so cell formatting/style
# indices remain unchanged from the original workbook.
$dRange.Style = "Normal"